$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Januari 2019")

# Replace personal/template data with placeholders (keep company & manager as-is)
$ws.Range("B2").Value = "<user>"
$ws.Range("B4").Value = "<month>"
$ws.Range("B5").Value = "<team>"
$ws.Range("B7").Value = "<projno>"

# Clear out sample daily hour entries
$ws.Range("B14").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("B28").ClearContents()

# Turn the monthly total row into a real formula instead of a static number
$ws.Range("B39").Formula = "=SUM(B8:B38)"

# Match the saved cursor/selection position recorded in the workbook
$ws.Range("O8").Select()
